$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column F: "time_taken" metadata, matching the header style of E1
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats

# time_taken values for data rows 2-10
$ws.Range("F2").Value = "2021-10-05 13:41:15.936612"
$ws.Range("F3").Value = "2021-10-05 13:41:15.936625"
$ws.Range("F4").Value = "2021-10-05 13:41:15.936629"
$ws.Range("F5").Value = "2021-10-05 13:41:15.936632"
$ws.Range("F6").Value = "2021-10-05 13:41:15.936635"
$ws.Range("F7").Value = "2021-10-05 13:41:15.936639"
$ws.Range("F8").Value = "2021-10-05 13:41:15.936642"
$ws.Range("F9").Value = "2021-10-05 13:41:15.936645"
$ws.Range("F10").Value = "2021-10-05 13:41:15.936648"
